$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header row values -> lowercase
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "path"
$ws.Range("C1").Value = "alias"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "card."
$ws.Range("F1").Value = "stereotype"
$ws.Range("G1").Value = "id"
$ws.Range("H1").Value = "definition"
$ws.Range("I1").Value = "definitioncode"

# Rows 6-13: fix column A (element name) and column B (path) values
$ws.Range("A6").Value = "ProductDescription"
$ws.Range("B6").Value = "MedicalDevice.ProductDescription"

$ws.Range("A7").Value = "AnatomicalLocation"
$ws.Range("B7").Value = "MedicalDevice.AnatomicalLocation"

$ws.Range("A8").Value = "Indication::Problem"
$ws.Range("B8").Value = "MedicalDevice.Indication::Problem"

$ws.Range("A9").Value = "StartDate"
$ws.Range("B9").Value = "MedicalDevice.StartDate"

$ws.Range("A10").Value = "EindDatum"
$ws.Range("B10").Value = "MedicalDevice.EindDatum"

$ws.Range("A11").Value = "Comment"
$ws.Range("B11").Value = "MedicalDevice.Comment"

$ws.Range("A12").Value = "Location::HealthcareProvider"
$ws.Range("B12").Value = "MedicalDevice.Location::HealthcareProvider"

$ws.Range("A13").Value = "HealthProfessional"
$ws.Range("B13").Value = "MedicalDevice.HealthProfessional"
